$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row "Création de la base de données": C, C, C -> R, R, R (Matthias's "AR" stays) ---
$t.Cell(5, 2).Range.Text = "R"
$t.Cell(5, 3).Range.Text = "R"
$t.Cell(5, 4).Range.Text = "R"

# --- Row "JSX" task renamed to "Style" ---
$t.Cell(8, 1).Range.Text = "Style"

# --- Row "CSS" task renamed to "Requêtes HTTP" ---
$t.Cell(9, 1).Range.Text = "Requêtes HTTP"

# --- Row "Algorithmique": Lucas R -> A + R (two runs, same formatting), Matthias AR -> R ---
$algoCell = $t.Cell(12, 2)
$algoRange = $algoCell.Range
$frag = '<w:r xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>A</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>R</w:t></w:r>'
$algoRange.InsertXML($frag)
$t.Cell(12, 5).Range.Text = "R"

# --- Row "Réception/Émission en JSON" renamed to "Routes" ---
$t.Cell(14, 1).Range.Text = "Routes"

Write-Output "changes applied"
